$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.338.20"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "3.150.22"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "591.48"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").Value = "3.141.95"
$ws.Range("E8").Value = "  +0.66%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.60%  "

$ws.Range("E10").Value = "  +0.43%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.85"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.80%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.27"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.73%  "

$ws.Range("D15").Value = "3.676.85"
$ws.Range("E15").Value = "  +0.93%  "

$ws.Range("E16").Value = "  -1.17%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.24"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.13%  "

$ws.Range("D18").Value = "64.156.20"
$ws.Range("E18").Value = "  +1.11%  "

$ws.Range("D19").Value = "3.152.70"
$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "466.87"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.72%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.733"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "81.26"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.02%  "

$ws.Range("E27").Value = "  +0.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.79"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.81%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.27"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.35%  "

$ws.Range("E32").Value = "  +0.16%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.50"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.57%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.110"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.02%  "

$ws.Range("D35").Value = "0.0₃0833"
$ws.Range("E35").Value = "  -4.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.23%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.02%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.31"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.23"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.96%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "465.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "51.53"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.97%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.297"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.75%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0375"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.83%  "

$ws.Range("D45").Value = "2.925.93"
$ws.Range("E45").Value = "  +0.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.81"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.108"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.88%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "129.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.67%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("E50").Value = "  +3.40%  "

$ws.Range("E51").Value = "  -0.42%  "
